$wb = $excel.ActiveWorkbook

# The workbook has two sheets that duplicate the same event rows:
#  - "展览" (sheet1)
#  - "全部类型" (sheet4)
# Both need the "想去人数" (column F) counts bumped for the same four rows.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 11775
    $ws.Range("F5").Value = 917
    $ws.Range("F17").Value = 1396
    $ws.Range("F19").Value = 911
}
